$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value2 = 93697.2651603072
$ws.Range("C2").Value2 = 86295.2695858756
$ws.Range("D2").Value2 = 82386.848118495
$ws.Range("E2").Value2 = 101165.614062139
$ws.Range("F2").Value2 = 104368.936761022
$ws.Range("I2").Value2 = 26015.2651603072

$ws.Range("B3").Value2 = 81377.2549593673
$ws.Range("C3").Value2 = 73034.8947051032
$ws.Range("D3").Value2 = 67869.9414167318
$ws.Range("E3").Value2 = 89895.4661580315
$ws.Range("F3").Value2 = 95877.2015995494
$ws.Range("I3").Value2 = 36444.2549593673

$ws.Range("B4").Value2 = 99663.735184421
$ws.Range("C4").Value2 = 91584.0911070704
$ws.Range("D4").Value2 = 87235.3050935919
$ws.Range("E4").Value2 = 108286.189852255
$ws.Range("F4").Value2 = 111699.631292856
$ws.Range("I4").Value2 = 26236.735184421

$ws.Range("B5").Value2 = 93810.1672091607
$ws.Range("C5").Value2 = 84714.2303745136
$ws.Range("D5").Value2 = 77678.1315738597
$ws.Range("E5").Value2 = 104366.643576462
$ws.Range("F5").Value2 = 109124.569499084
$ws.Range("I5").Value2 = 8126.16720916075

$ws.Range("B6").Value2 = 91147.5374579115
$ws.Range("C6").Value2 = 81552.7691434707
$ws.Range("D6").Value2 = 75474.4763119924
$ws.Range("E6").Value2 = 101686.749493398
$ws.Range("F6").Value2 = 107788.395387659
$ws.Range("I6").Value2 = 7762.53745791155

$ws.Range("B7").Value2 = 87505.0477078806
$ws.Range("C7").Value2 = 75899.9659327145
$ws.Range("D7").Value2 = 69104.7510090675
$ws.Range("E7").Value2 = 99236.5238395072
$ws.Range("F7").Value2 = 105024.576809057
$ws.Range("I7").Value2 = 2553.04770788057

$ws.Range("B8").Value2 = 85160.3861940999
$ws.Range("C8").Value2 = 73251.2736705855
$ws.Range("D8").Value2 = 66411.6356203747
$ws.Range("E8").Value2 = 96135.7473181393
$ws.Range("F8").Value2 = 103715.112791038
$ws.Range("I8").Value2 = 8737.38619409995

$ws.Range("B9").Value2 = 81941.051850653
$ws.Range("C9").Value2 = 68584.2504679944
$ws.Range("D9").Value2 = 60907.3833769565
$ws.Range("E9").Value2 = 95367.8435419574
$ws.Range("F9").Value2 = 101075.53398921
$ws.Range("I9").Value2 = -1159.94814934698

$ws.Range("B10").Value2 = 77473.6827770587
$ws.Range("C10").Value2 = 62565.1404407029
$ws.Range("D10").Value2 = 55050.0128236827
$ws.Range("E10").Value2 = 92058.3872891522
$ws.Range("F10").Value2 = 100013.303556282
$ws.Range("I10").Value2 = 2064.68277705867

$ws.Range("B11").Value2 = 74391.7212347732
$ws.Range("C11").Value2 = 59021.5858695057
$ws.Range("D11").Value2 = 49440.0709255898
$ws.Range("E11").Value2 = 89081.3625651318
$ws.Range("F11").Value2 = 96395.6512992942
$ws.Range("I11").Value2 = 6548.7212347732

$ws.Range("B12").Value2 = 71600.4241410029
$ws.Range("C12").Value2 = 55846.9031997573
$ws.Range("D12").Value2 = 47155.2162392888
$ws.Range("E12").Value2 = 85626.4098931422
$ws.Range("F12").Value2 = 93347.1196862182
$ws.Range("I12").Value2 = 1960.42414100292

$ws.Range("B13").Value2 = 69076.5048106195
$ws.Range("C13").Value2 = 51805.089877229
$ws.Range("D13").Value2 = 41615.6849803005
$ws.Range("E13").Value2 = 85454.7968545395
$ws.Range("F13").Value2 = 91114.7980643164
$ws.Range("I13").Value2 = 4979.50481061952

$ws.Range("B14").Value2 = 91769.5363592395
$ws.Range("C14").Value2 = 74223.5147925542
$ws.Range("D14").Value2 = 64074.1388191765
$ws.Range("E14").Value2 = 109065.842061821
$ws.Range("F14").Value2 = 117423.564707125
$ws.Range("I14").Value2 = 26956.5363592395

$ws.Range("B15").Value2 = 79779.3556617333
$ws.Range("C15").Value2 = 61131.7259826951
$ws.Range("D15").Value2 = 52690.3477543721
$ws.Range("E15").Value2 = 96593.7330100993
$ws.Range("F15").Value2 = 106973.00587313
$ws.Range("I15").Value2 = 24354.3556617333

$ws.Range("B16").Value2 = 98137.1461899836
$ws.Range("C16").Value2 = 80022.5405824678
$ws.Range("D16").Value2 = 67042.8751685108
$ws.Range("E16").Value2 = 115109.62089298
$ws.Range("F16").Value2 = 124847.969807371
$ws.Range("I16").Value2 = 17334.1461899836

$ws.Range("B17").Value2 = 92638.463776511
$ws.Range("C17").Value2 = 72838.9448037355
$ws.Range("D17").Value2 = 60438.7976921127
$ws.Range("E17").Value2 = 110858.344307543
$ws.Range("F17").Value2 = 121313.374481702
$ws.Range("I17").Value2 = 12090.463776511

$ws.Range("B18").Value2 = 89478.2399346006
$ws.Range("C18").Value2 = 69700.2727437459
$ws.Range("D18").Value2 = 59399.4967794055
$ws.Range("E18").Value2 = 108964.808400609
$ws.Range("F18").Value2 = 121099.10440871
$ws.Range("I18").Value2 = 14235.2399346006

$ws.Range("B19").Value2 = 85762.555098738
$ws.Range("C19").Value2 = 65415.0534934471
$ws.Range("D19").Value2 = 54749.8651722092
$ws.Range("E19").Value2 = 106273.224604777
$ws.Range("F19").Value2 = 118308.275708742
$ws.Range("I19").Value2 = 11878.555098738

$ws.Range("B20").Value2 = 83285.0403004039
$ws.Range("C20").Value2 = 63505.16476858
$ws.Range("D20").Value2 = 48397.9991052699
$ws.Range("E20").Value2 = 103434.047709919
$ws.Range("F20").Value2 = 117672.090787915
$ws.Range("I20").Value2 = 6637.04030040387

$ws.Range("B21").Value2 = 80373.1093157834
$ws.Range("C21").Value2 = 58300.639041366
$ws.Range("D21").Value2 = 42135.4237414375
$ws.Range("E21").Value2 = 102870.677736677
$ws.Range("F21").Value2 = 113715.498668038
$ws.Range("I21").Value2 = 12407.1093157834

$ws.Range("B22").Value2 = 76130.4155476807
$ws.Range("C22").Value2 = 54403.76940199
$ws.Range("D22").Value2 = 33974.2956101364
$ws.Range("E22").Value2 = 97683.9822561114
$ws.Range("F22").Value2 = 111114.182814459
$ws.Range("I22").Value2 = 8318.41554768071

$ws.Range("B23").Value2 = 73042.9692172477
$ws.Range("C23").Value2 = 49446.447314652
$ws.Range("D23").Value2 = 34315.5549174296
$ws.Range("E23").Value2 = 96049.8544086778
$ws.Range("F23").Value2 = 110142.06310894
$ws.Range("I23").Value2 = 11651.9692172477

$ws.Range("B24").Value2 = 70291.7678990695
$ws.Range("C24").Value2 = 46717.5955634672
$ws.Range("D24").Value2 = 34491.9253428723
$ws.Range("E24").Value2 = 93526.5051967359
$ws.Range("F24").Value2 = 105932.469970768
$ws.Range("I24").Value2 = 8538.76789906953

$ws.Range("B25").Value2 = 67713.5703077417
$ws.Range("C25").Value2 = 42867.9518203189
$ws.Range("D25").Value2 = 32225.1800955643
$ws.Range("E25").Value2 = 92369.7172924127
$ws.Range("F25").Value2 = 105882.380893609
$ws.Range("I25").Value2 = 5925.5703077417

$ws.Range("B26").Value2 = 90408.6807052656
$ws.Range("C26").Value2 = 65609.4273309539
$ws.Range("D26").Value2 = 46822.7786045865
$ws.Range("E26").Value2 = 116176.603736778
$ws.Range("F26").Value2 = 126482.039548605
$ws.Range("I26").Value2 = 28711.6807052656

$ws.Range("B27").Value2 = 78599.4186776249
$ws.Range("C27").Value2 = 52145.7433841234
$ws.Range("D27").Value2 = 36646.2185060523
$ws.Range("E27").Value2 = 103958.08916066
$ws.Range("F27").Value2 = 116390.908353683
$ws.Range("I27").Value2 = 26003.4186776249

$ws.Range("B28").Value2 = 97253.8889243581
$ws.Range("C28").Value2 = 71602.3220136479
$ws.Range("D28").Value2 = 55620.1772945534
$ws.Range("E28").Value2 = 121938.151087739
$ws.Range("F28").Value2 = 135415.432575173
$ws.Range("I28").Value2 = 24143.8889243581

$ws.Range("B29").Value2 = 91776.9576062413
$ws.Range("C29").Value2 = 65931.5121417484
$ws.Range("D29").Value2 = 44908.2290072984
$ws.Range("E29").Value2 = 115681.393232848
$ws.Range("F29").Value2 = 131596.924416579
$ws.Range("I29").Value2 = 30591.9576062413

$ws.Range("B30").Value2 = 88082.9425920034
$ws.Range("C30").Value2 = 60711.1125331854
$ws.Range("D30").Value2 = 41788.4484945635
$ws.Range("E30").Value2 = 113567.613115531
$ws.Range("F30").Value2 = 131373.347987105
$ws.Range("I30").Value2 = 24492.9425920034

$ws.Range("B31").Value2 = 84818.1121163959
$ws.Range("C31").Value2 = 56974.7351179704
$ws.Range("D31").Value2 = 40286.5254321612
$ws.Range("E31").Value2 = 111015.737869733
$ws.Range("F31").Value2 = 125078.356899127
$ws.Range("I31").Value2 = 16917.1121163959

$ws.Range("B32").Value2 = 82563.2233071128
$ws.Range("C32").Value2 = 54224.734838107
$ws.Range("D32").Value2 = 37105.4619654902
$ws.Range("E32").Value2 = 110799.944554892
$ws.Range("F32").Value2 = 125462.908927257
$ws.Range("I32").Value2 = 11141.2233071128

$ws.Range("B33").Value2 = 79607.7557549771
$ws.Range("C33").Value2 = 51878.2370241907
$ws.Range("D33").Value2 = 34177.2217371328
$ws.Range("E33").Value2 = 106551.137959744
$ws.Range("F33").Value2 = 122891.571466379
$ws.Range("I33").Value2 = 10588.7557549771

$ws.Range("B34").Value2 = 75082.3764572099
$ws.Range("C34").Value2 = 44222.2206530664
$ws.Range("D34").Value2 = 28300.0410984712
$ws.Range("E34").Value2 = 104503.109026289
$ws.Range("F34").Value2 = 120660.079318795
$ws.Range("I34").Value2 = 16444.3764572099

$ws.Range("B35").Value2 = 72333.8593665555
$ws.Range("C35").Value2 = 40714.8322744468
$ws.Range("D35").Value2 = 26501.7935734671
$ws.Range("E35").Value2 = 101739.672142452
$ws.Range("F35").Value2 = 119895.502132153
$ws.Range("I35").Value2 = 21208.8593665555

$ws.Range("B36").Value2 = 69987.2202953226
$ws.Range("C36").Value2 = 39712.0019655725
$ws.Range("D36").Value2 = 23782.0209564259
$ws.Range("E36").Value2 = 99573.6542300597
$ws.Range("F36").Value2 = 114910.667965309
$ws.Range("I36").Value2 = 21635.2202953226

$ws.Range("B37").Value2 = 67338.0814954225
$ws.Range("C37").Value2 = 34504.1226858035
$ws.Range("D37").Value2 = 21507.7784524371
$ws.Range("E37").Value2 = 97081.7384077673
$ws.Range("F37").Value2 = 116361.723447662
$ws.Range("I37").Value2 = 33387.0814954225

$ws.Range("B38").Value2 = 89966.1808380493
$ws.Range("C38").Value2 = 57280.1032027045
$ws.Range("D38").Value2 = 41504.1038846977
$ws.Range("E38").Value2 = 121328.616067055
$ws.Range("F38").Value2 = 135900.000875582
$ws.Range("I38").Value2 = 36236.1808380493

$ws.Range("B39").Value2 = 77415.1980253574
$ws.Range("C39").Value2 = 43170.8255278554
$ws.Range("D39").Value2 = 27375.812300636
$ws.Range("E39").Value2 = 108751.252098411
$ws.Range("F39").Value2 = 127047.854164887
$ws.Range("I39").Value2 = 5574.1980253574

$ws.Range("B40").Value2 = 95911.2167193544
$ws.Range("C40").Value2 = 60777.8014929914
$ws.Range("D40").Value2 = 45286.4152194862
$ws.Range("E40").Value2 = 129794.167619622
$ws.Range("F40").Value2 = 147811.583374455
$ws.Range("I40").Value2 = 19580.2167193544
